$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("General Settings")

# The workbook links to an external add-in workbook (FixedIncome.xla) via
# xl/externalLinks - cell D8 holds a formula that calls into it
# ([1]!qlSerializationPath(Trigger)). The commit removes that link and
# replaces the formula with a plain literal path.

# 1) Stamp the new literal value into D8 as a value (not a formula) while
#    the external link is still intact, using a value-only paste so the
#    cell's existing number format / style survives untouched (a plain
#    .Value assignment here would make Excel silently renormalize the
#    cell's quote-prefixed style).
$ws.Range("Z1").Value = "C:\Users\erik\junk\"
$ws.Range("Z1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("Z1").Clear()

# 2) Now that no formula references the external workbook any more, break
#    the link so Excel drops the externalReferences/externalLinks parts
#    entirely.
$wb.BreakLink("/WorkGroup/IMI_Workbooks/Production/QLXL_R01030x/framework/addin/FixedIncome.xla", 1)
